$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A72").Value = "2025/12/05 18:00"
$ws.Range("B72").Value = "11,730位本"
$ws.Range("C72").Value = "37位 広告・宣伝 (本)"
$ws.Range("D72").Value = "48位商業デザイン"
$ws.Range("E72").Value = "756位ビジネス実用本"
$ws.Range("F72").Value = "-"
$ws.Range("G72").Value = "-"
